$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory (column H) labels for specific rows
$ws.Range("H7").Value  = "line graph(s)"
$ws.Range("H9").Value  = "photo(s)"
$ws.Range("H10").Value = "data display"
$ws.Range("H11").Value = "data display"
$ws.Range("H15").Value = "bar chart(s)"
$ws.Range("H16").Value = "line graph(s)"
$ws.Range("H23").Value = "photo(s)"

# Remove the entire "is_viewed" column (column I)
$ws.Columns.Item(9).Delete()
